# ChainableWaits + step numbers + colors + rework driver restart + PageElement & Page cucumber param
#
# The "Result" column (F) of the NoraUi-blog demo sheet used to show the raw
# localized messages produced by the step results ("Succès", "Échec : ...",
# "Élément ignoré ..."). This updates the sample data so the column instead
# shows the new generic chainable-wait outcome labels ("Success" / "Fail: ...")
# with their status colors (green = success, red = failure).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$successText = "Success"
$failText = "Fail: anonymous is prohibited in demo blog!!"

# Excel standard palette: indexed 17 = green, indexed 10 = red.
$greenColor = 32768   # RGB(0,128,0)
$redColor   = 255     # RGB(255,0,0)

$successCells = @("F2", "F3", "F4", "F7")
foreach ($addr in $successCells) {
    $cell = $ws.Range($addr)
    $cell.Value = $successText
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Color = $greenColor
}

$failCell = $ws.Range("F6")
$failCell.Value = $failText
$failCell.Font.Name = "Calibri"
$failCell.Font.Size = 11
$failCell.Font.Color = $redColor
